# Natmi following Dr Hou advice
# Update Ligand/Receptor-expressing cell counts (3 instead of 1) and
# recompute dependent expression/specificity values for rows 2-17.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 157.7984646666667
$ws.Range("H2").Value = 473.395394
$ws.Range("I2").Value = 0.341075365555871
$ws.Range("J2").Value = 0.3410753655558709
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.854571666666667
$ws.Range("N2").Value = 5.563715
$ws.Range("O2").Value = 0.01651371646154392
$ws.Range("P2").Value = 0.01651371646154392
$ws.Range("Q2").Value = 292.6485616143011
$ws.Range("R2").Value = 2633.83705452871
$ws.Range("S2").Value = 0.005632421878807096
$ws.Range("T2").Value = 0.005632421878807097
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 157.7984646666667
$ws.Range("H3").Value = 473.395394
$ws.Range("I3").Value = 0.341075365555871
$ws.Range("J3").Value = 0.3410753655558709
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 92.91372433333333
$ws.Range("N3").Value = 278.741173
$ws.Range("O3").Value = 0.8273343794712995
$ws.Range("P3").Value = 0.8273343794712996
$ws.Range("Q3").Value = 14661.64304626191
$ws.Range("R3").Value = 131954.7874163572
$ws.Range("S3").Value = 0.2821833759151132
$ws.Range("T3").Value = 0.2821833759151132
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 157.7984646666667
$ws.Range("H4").Value = 473.395394
$ws.Range("I4").Value = 0.341075365555871
$ws.Range("J4").Value = 0.3410753655558709
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.162136
$ws.Range("N4").Value = 0.4864080000000001
$ws.Range("O4").Value = 0.001443712303133186
$ws.Range("P4").Value = 0.001443712303133187
$ws.Range("Q4").Value = 25.58481186719467
$ws.Range("R4").Value = 230.263306804752
$ws.Range("S4").Value = 0.0004924147015486599
$ws.Range("T4").Value = 0.0004924147015486599
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 157.7984646666667
$ws.Range("H5").Value = 473.395394
$ws.Range("I5").Value = 0.341075365555871
$ws.Range("J5").Value = 0.3410753655558709
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.37449166666667
$ws.Range("N5").Value = 52.123475
$ws.Range("O5").Value = 0.1547081917640233
$ws.Range("P5").Value = 0.1547081917640233
$ws.Range("Q5").Value = 2741.668109363794
$ws.Range("R5").Value = 24675.01298427415
$ws.Range("S5").Value = 0.05276715306040204
$ws.Range("T5").Value = 0.05276715306040203
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 127.5109433333333
$ws.Range("H6").Value = 382.53283
$ws.Range("I6").Value = 0.2756100428585324
$ws.Range("J6").Value = 0.2756100428585324
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.854571666666667
$ws.Range("N6").Value = 5.563715
$ws.Range("O6").Value = 0.01651371646154392
$ws.Range("P6").Value = 0.01651371646154392
$ws.Range("Q6").Value = 236.4781826959389
$ws.Range("R6").Value = 2128.30364426345
$ws.Range("S6").Value = 0.004551346101719772
$ws.Range("T6").Value = 0.004551346101719773
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 127.5109433333333
$ws.Range("H7").Value = 382.53283
$ws.Range("I7").Value = 0.2756100428585324
$ws.Range("J7").Value = 0.2756100428585324
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 92.91372433333333
$ws.Range("N7").Value = 278.741173
$ws.Range("O7").Value = 0.8273343794712995
$ws.Range("P7").Value = 0.8273343794712996
$ws.Range("Q7").Value = 11847.51663835662
$ws.Range("R7").Value = 106627.6497452096
$ws.Range("S7").Value = 0.2280216637844222
$ws.Range("T7").Value = 0.2280216637844222
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 127.5109433333333
$ws.Range("H8").Value = 382.53283
$ws.Range("I8").Value = 0.2756100428585324
$ws.Range("J8").Value = 0.2756100428585324
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.162136
$ws.Range("N8").Value = 0.4864080000000001
$ws.Range("O8").Value = 0.001443712303133186
$ws.Range("P8").Value = 0.001443712303133187
$ws.Range("Q8").Value = 20.67411430829334
$ws.Range("R8").Value = 186.06702877464
$ws.Range("S8").Value = 0.000397901609741928
$ws.Range("T8").Value = 0.000397901609741928
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 127.5109433333333
$ws.Range("H9").Value = 382.53283
$ws.Range("I9").Value = 0.2756100428585324
$ws.Range("J9").Value = 0.2756100428585324
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.37449166666667
$ws.Range("N9").Value = 52.123475
$ws.Range("O9").Value = 0.1547081917640233
$ws.Range("P9").Value = 0.1547081917640233
$ws.Range("Q9").Value = 2215.437822353806
$ws.Range("R9").Value = 19938.94040118425
$ws.Range("S9").Value = 0.04263913136264852
$ws.Range("T9").Value = 0.04263913136264852
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 86.127454
$ws.Range("H10").Value = 258.382362
$ws.Range("I10").Value = 0.1861612083457225
$ws.Range("J10").Value = 0.1861612083457225
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.854571666666667
$ws.Range("N10").Value = 5.563715
$ws.Range("O10").Value = 0.01651371646154392
$ws.Range("P10").Value = 0.01651371646154392
$ws.Range("Q10").Value = 159.7295359105367
$ws.Range("R10").Value = 1437.56582319483
$ws.Range("S10").Value = 0.003074213410759664
$ws.Range("T10").Value = 0.003074213410759665
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 86.127454
$ws.Range("H11").Value = 258.382362
$ws.Range("I11").Value = 0.1861612083457225
$ws.Range("J11").Value = 0.1861612083457225
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 92.91372433333333
$ws.Range("N11").Value = 278.741173
$ws.Range("O11").Value = 0.8273343794712995
$ws.Range("P11").Value = 0.8273343794712996
$ws.Range("Q11").Value = 8002.422518487848
$ws.Range("R11").Value = 72021.80266639062
$ws.Range("S11").Value = 0.1540175677883356
$ws.Range("T11").Value = 0.1540175677883356
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 86.127454
$ws.Range("H12").Value = 258.382362
$ws.Range("I12").Value = 0.1861612083457225
$ws.Range("J12").Value = 0.1861612083457225
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.162136
$ws.Range("N12").Value = 0.4864080000000001
$ws.Range("O12").Value = 0.001443712303133186
$ws.Range("P12").Value = 0.001443712303133187
$ws.Range("Q12").Value = 13.964360881744
$ws.Range("R12").Value = 125.679247935696
$ws.Range("S12").Value = 0.0002687632268548599
$ws.Range("T12").Value = 0.00026876322685486
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 86.127454
$ws.Range("H13").Value = 258.382362
$ws.Range("I13").Value = 0.1861612083457225
$ws.Range("J13").Value = 0.1861612083457225
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.37449166666667
$ws.Range("N13").Value = 52.123475
$ws.Range("O13").Value = 0.1547081917640233
$ws.Range("P13").Value = 0.1547081917640233
$ws.Range("Q13").Value = 1496.420731794217
$ws.Range("R13").Value = 13467.78658614795
$ws.Range("S13").Value = 0.02880066391977233
$ws.Range("T13").Value = 0.02880066391977233
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 91.212982
$ws.Range("H14").Value = 273.638946
$ws.Range("I14").Value = 0.1971533832398742
$ws.Range("J14").Value = 0.1971533832398741
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.854571666666667
$ws.Range("N14").Value = 5.563715
$ws.Range("O14").Value = 0.01651371646154392
$ws.Range("P14").Value = 0.01651371646154392
$ws.Range("Q14").Value = 169.1610120493767
$ws.Range("R14").Value = 1522.44910844439
$ws.Range("S14").Value = 0.003255735070257387
$ws.Range("T14").Value = 0.003255735070257387
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 91.212982
$ws.Range("H15").Value = 273.638946
$ws.Range("I15").Value = 0.1971533832398742
$ws.Range("J15").Value = 0.1971533832398741
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 92.91372433333333
$ws.Range("N15").Value = 278.741173
$ws.Range("O15").Value = 0.8273343794712995
$ws.Range("P15").Value = 0.8273343794712996
$ws.Range("Q15").Value = 8474.937865169295
$ws.Range("R15").Value = 76274.44078652366
$ws.Range("S15").Value = 0.1631117719834286
$ws.Range("T15").Value = 0.1631117719834286
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 91.212982
$ws.Range("H16").Value = 273.638946
$ws.Range("I16").Value = 0.1971533832398742
$ws.Range("J16").Value = 0.1971533832398741
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.162136
$ws.Range("N16").Value = 0.4864080000000001
$ws.Range("O16").Value = 0.001443712303133186
$ws.Range("P16").Value = 0.001443712303133187
$ws.Range("Q16").Value = 14.788908049552
$ws.Range("R16").Value = 133.100172445968
$ws.Range("S16").Value = 0.0002846327649877385
$ws.Range("T16").Value = 0.0002846327649877385
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 91.212982
$ws.Range("H17").Value = 273.638946
$ws.Range("I17").Value = 0.1971533832398742
$ws.Range("J17").Value = 0.1971533832398741
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.37449166666667
$ws.Range("N17").Value = 52.123475
$ws.Range("O17").Value = 0.1547081917640233
$ws.Range("P17").Value = 0.1547081917640233
$ws.Range("Q17").Value = 1584.779195650817
$ws.Range("R17").Value = 14263.01276085735
$ws.Range("S17").Value = 0.03050124342120043
$ws.Range("T17").Value = 0.03050124342120043
